$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: Row content swaps / rotations.
# Several match rows were re-ordered (their whole data row, columns B..AB,
# moved to a different row position) while the running index in column A
# stayed tied to the row position.
# ---------------------------------------------------------------------------

function Swap-Rows($rowA, $rowB) {
    $rangeA = "B$rowA`:AB$rowA"
    $rangeB = "B$rowB`:AB$rowB"
    $valA = $ws.Range($rangeA).Value2
    $valB = $ws.Range($rangeB).Value2
    $ws.Range($rangeA).Value2 = $valB
    $ws.Range($rangeB).Value2 = $valA
}

# Simple two-row swaps
Swap-Rows 6 7
Swap-Rows 24 25
Swap-Rows 32 33
Swap-Rows 40 41
Swap-Rows 77 78
Swap-Rows 164 165
Swap-Rows 170 171
Swap-Rows 197 199

# Three-row rotations: new(61) = old(63); new(62) = old(61); new(63) = old(62)
$r61 = $ws.Range("B61:AB61").Value2
$r62 = $ws.Range("B62:AB62").Value2
$r63 = $ws.Range("B63:AB63").Value2
$ws.Range("B61:AB61").Value2 = $r63
$ws.Range("B62:AB62").Value2 = $r61
$ws.Range("B63:AB63").Value2 = $r62

# new(204) = old(206); new(205) = old(204); new(206) = old(205)
$r204 = $ws.Range("B204:AB204").Value2
$r205 = $ws.Range("B205:AB205").Value2
$r206 = $ws.Range("B206:AB206").Value2
$ws.Range("B204:AB204").Value2 = $r206
$ws.Range("B205:AB205").Value2 = $r204
$ws.Range("B206:AB206").Value2 = $r205

# ---------------------------------------------------------------------------
# Part 2: Updated odds figures for upcoming fixtures (rows 208-215).
# ---------------------------------------------------------------------------

$ws.Range("M208").Value2  = 2.05
$ws.Range("Q208").Value2  = 2.08
$ws.Range("R208").Value2  = 1.82

$ws.Range("M209").Value2  = 2.375
$ws.Range("N209").Value2  = 3.3
$ws.Range("O209").Value2  = 3

$ws.Range("M210").Value2  = 2.25
$ws.Range("N210").Value2  = 4
$ws.Range("O210").Value2  = 2.8
$ws.Range("P210").Value2  = -0.25
$ws.Range("Q210").Value2  = 2.05
$ws.Range("R210").Value2  = 1.85
$ws.Range("T210").Value2  = 2.025
$ws.Range("U210").Value2  = 1.825

$ws.Range("M211").Value2  = 3
$ws.Range("O211").Value2  = 2.375
$ws.Range("Q211").Value2  = 1.86
$ws.Range("R211").Value2  = 2.04

$ws.Range("Q213").Value2  = 1.93
$ws.Range("R213").Value2  = 1.97

$ws.Range("Q214").Value2  = 1.88
$ws.Range("R214").Value2  = 2.02

$ws.Range("Q215").Value2  = 2.11
$ws.Range("R215").Value2  = 1.79
